$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.331.44"
$ws.Range("E2").Value = "  -4.75%  "
$ws.Range("D3").Value = "3.273.15"
$ws.Range("E3").Value = "  -7.17%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -11.62%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "3.266.90"
$ws.Range("E8").Value = "  -7.25%  "
$ws.Range("E9").Value = "  -10.80%  "
$ws.Range("E10").Value = "  -13.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.68"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.513"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -12.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.61"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -16.42%  "
$ws.Range("E14").Value = "  -11.11%  "
$ws.Range("D15").Value = "3.797.19"
$ws.Range("E15").Value = "  -7.15%  "
$ws.Range("D16").Value = "67.351.40"
$ws.Range("E16").Value = "  -4.87%  "
$ws.Range("D17").Value = "3.269.62"
$ws.Range("E17").Value = "  -7.18%  "
$ws.Range("E18").Value = "  -13.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "535.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -11.85%  "
$ws.Range("E20").Value = "  -6.11%  "
$ws.Range("E21").Value = "  -14.31%  "
$ws.Range("E22").Value = "  -13.20%  "
$ws.Range("E23").Value = "  -12.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -11.82%  "
$ws.Range("E25").Value = "  -12.29%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  -12.42%  "
$ws.Range("E28").Value = "  -10.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "29.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -12.39%  "
$ws.Range("E30").Value = "  -16.01%  "
$ws.Range("E31").Value = "  -10.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -11.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "544.82"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -11.48%  "
$ws.Range("E34").Value = "  -18.08%  "
$ws.Range("E35").Value = "  -15.02%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0458"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "53.50"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.64%  "
$ws.Range("E39").Value = "  -13.08%  "
$ws.Range("E40").Value = "  -15.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.129"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -9.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -18.87%  "
$ws.Range("D43").Value = "2.948.28"
$ws.Range("E43").Value = "  -11.76%  "
$ws.Range("D44").Value = "0.0₃0596"
$ws.Range("E44").Value = "  -17.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.265"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -14.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "26.93"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -15.31%  "
$ws.Range("E47").Value = "  -12.98%  "
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("E49").Value = "  -19.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "126.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.56%  "
$ws.Range("E51").Value = "  -12.32%  "
